$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the whole data records currently stored in row 19
# ("Gultoppig fingersvamp") and row 20 ("Knärot") - the two observations
# change places in the sheet while every other row stays untouched.
#
# Only the columns whose values actually differ between the two rows are
# touched (A, B, D, E, F, G, H, I, J, Q, R, Z, AB); this avoids Excel
# re-interpreting/auto-converting any of the other, unrelated columns
# (e.g. the date columns) that are left alone by the real edit.

$numericTextCols = @("I")   # columns that hold numeric-looking text values
$cols = @("A","B","D","E","F","G","H","I","J","Q","R","Z","AB")

foreach ($col in $cols) {
    $cell19 = $ws.Range("$col`19")
    $cell20 = $ws.Range("$col`20")

    $v19 = $cell19.Value()
    $v20 = $cell20.Value()

    if ($numericTextCols -contains $col) {
        # Force a text cell so numeric-looking values ("20", "25", ...)
        # are not silently turned into real numbers by Excel.
        $cell19.NumberFormat = "@"
        $cell20.NumberFormat = "@"
    }

    $cell19.Value = $v20
    $cell20.Value = $v19
}

# Column L on row 20 is an empty placeholder cell before the edit; once
# the records swap places it no longer belongs to row 20 (it moves over
# to row 19, which has no value there either way).
$ws.Range("L20").ClearContents()
